$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.32672348343714
$ws.Range("C2").Value = 7.032547231609001
$ws.Range("D2").Value = 13.23667059365011
$ws.Range("E2").Value = 13.30050451092633
$ws.Range("G2").Value = 3.691422788458469
$ws.Range("I2").Value = 27.80625263008561
$ws.Range("J2").Value = 8.124831417009732
$ws.Range("L2").Value = 13.02087881682436
$ws.Range("N2").Value = 18.87737904267741
$ws.Range("O2").Value = 30.01451601961928
$ws.Range("B3").Value = 19.94821412473954
$ws.Range("C3").Value = 6.636453385791423
$ws.Range("D3").Value = 13.24272221209651
$ws.Range("E3").Value = 13.3321159949392
$ws.Range("G3").Value = 3.693932781911682
$ws.Range("I3").Value = 27.90737420016027
$ws.Range("J3").Value = 8.130868781462437
$ws.Range("L3").Value = 13.01102782136025
$ws.Range("N3").Value = 18.93706305013791
$ws.Range("O3").Value = 30.07524951261751
$ws.Range("B4").Value = 19.71555049186747
$ws.Range("C4").Value = 6.379477741944648
$ws.Range("D4").Value = 13.24881455245389
$ws.Range("E4").Value = 13.35312613052129
$ws.Range("G4").Value = 3.695556367417525
$ws.Range("I4").Value = 27.97503563685443
$ws.Range("J4").Value = 8.134789824018144
$ws.Range("L4").Value = 13.00691320339064
$ws.Range("N4").Value = 18.97558595449494
$ws.Range("O4").Value = 30.11932099638203
$ws.Range("B5").Value = 19.62079891522044
$ws.Range("C5").Value = 6.271332363913692
$ws.Range("D5").Value = 13.25189528186098
$ws.Range("E5").Value = 13.36209077031533
$ws.Range("G5").Value = 3.69623878979129
$ws.Range("I5").Value = 28.00400688858281
$ws.Range("J5").Value = 8.136441661939408
$ws.Range("L5").Value = 13.00572476837068
$ws.Range("N5").Value = 18.99175752699884
$ws.Range("O5").Value = 30.13898152235226
$ws.Range("B6").Value = 19.60507286503971
$ws.Range("C6").Value = 6.253169282884449
$ws.Range("D6").Value = 13.25244296390695
$ws.Range("E6").Value = 13.36360368312696
$ws.Range("G6").Value = 3.696353363662585
$ws.Range("I6").Value = 28.00890194383681
$ws.Range("J6").Value = 8.136719213118008
$ws.Range("L6").Value = 13.00555697609624
$ws.Range("N6").Value = 18.99447142049362
$ws.Range("O6").Value = 30.14234874951336
$ws.Range("B7").Value = 19.7142722257012
$ws.Range("C7").Value = 6.378033069798478
$ws.Range("D7").Value = 13.24885367834457
$ws.Range("E7").Value = 13.3532453992219
$ws.Range("G7").Value = 3.69556548650346
$ws.Range("I7").Value = 27.97542069356686
$ws.Range("J7").Value = 8.134811882513899
$ws.Range("L7").Value = 13.00689519613726
$ws.Range("N7").Value = 18.9758021324261
$ws.Range("O7").Value = 30.11957926338869
$ws.Range("B8").Value = 20.19634139952786
$ws.Range("C8").Value = 6.898852237251421
$ws.Range("D8").Value = 13.23826424090642
$ws.Range("E8").Value = 13.31107217620527
$ws.Range("G8").Value = 3.692271164926877
$ws.Range("I8").Value = 27.83996127222872
$ws.Range("J8").Value = 8.126868759087829
$ws.Range("L8").Value = 13.01708207098611
$ws.Range("N8").Value = 18.89756934678229
$ws.Range("O8").Value = 30.03404722746899
$ws.Range("B9").Value = 21.13426148060372
$ws.Range("C9").Value = 7.8096820624539
$ws.Range("D9").Value = 13.23632971600388
$ws.Range("E9").Value = 13.24105403938091
$ws.Range("G9").Value = 3.68646202856756
$ws.Range("I9").Value = 27.61865360537207
$ws.Range("J9").Value = 8.112984105625477
$ws.Range("L9").Value = 13.05230738774878
$ws.Range("N9").Value = 18.75898770547266
$ws.Range("O9").Value = 29.92028913337122
$ws.Range("B10").Value = 21.81189067101674
$ws.Range("C10").Value = 8.410361496727079
$ws.Range("D10").Value = 13.24634271760607
$ws.Range("E10").Value = 13.19732168986459
$ws.Range("G10").Value = 3.682586598448946
$ws.Range("I10").Value = 27.48323646230982
$ws.Range("J10").Value = 8.103804919977934
$ws.Range("L10").Value = 13.08734303270064
$ws.Range("N10").Value = 18.66613277712204
$ws.Range("O10").Value = 29.86981787890616
$ws.Range("B11").Value = 22.11629485899123
$ws.Range("C11").Value = 8.668614926958782
$ws.Range("D11").Value = 13.25336642573982
$ws.Range("E11").Value = 13.17909614466858
$ws.Range("G11").Value = 3.680907879351342
$ws.Range("I11").Value = 27.42756569763419
$ws.Range("J11").Value = 8.099848946882764
$ws.Range("L11").Value = 13.10523505134565
$ws.Range("N11").Value = 18.62581935339424
$ws.Range("O11").Value = 29.85408069549418
$ws.Range("B12").Value = 22.23090315140133
$ws.Range("C12").Value = 8.764245309681543
$ws.Range("D12").Value = 13.2563794801036
$ws.Range("E12").Value = 13.17243417219633
$ws.Range("G12").Value = 3.680284234615009
$ws.Range("I12").Value = 27.40734022522021
$ws.Range("J12").Value = 8.098382360365228
$ws.Range("L12").Value = 13.11228793655317
$ws.Range("N12").Value = 18.61082952141224
$ws.Range("O12").Value = 29.8491619788117
$ws.Range("B13").Value = 22.20625136155716
$ws.Range("C13").Value = 8.74374595769237
$ws.Range("D13").Value = 13.25571487975613
$ws.Range("E13").Value = 13.17385829342209
$ws.Range("G13").Value = 3.680418012748644
$ws.Range("I13").Value = 27.41165803254964
$ws.Range("J13").Value = 8.098696819207714
$ws.Range("L13").Value = 13.11075668386867
$ws.Range("N13").Value = 18.61404559227958
$ws.Range("O13").Value = 29.85017500409027
$ws.Range("B14").Value = 22.12573764177214
$ws.Range("C14").Value = 8.676525911295329
$ws.Range("D14").Value = 13.25360724074346
$ws.Range("E14").Value = 13.17854325896553
$ws.Range("G14").Value = 3.680856330587086
$ws.Range("I14").Value = 27.42588456687635
$ws.Range("J14").Value = 8.099727660302483
$ws.Range("L14").Value = 13.1058097541189
$ws.Range("N14").Value = 18.62458060580906
$ws.Range("O14").Value = 29.85365516448904
$ws.Range("B15").Value = 22.07633119439561
$ws.Range("C15").Value = 8.635069620704227
$ws.Range("D15").Value = 13.25236221103472
$ws.Range("E15").Value = 13.18144413610018
$ws.Range("G15").Value = 3.68112638029353
$ws.Range("I15").Value = 27.43471027190859
$ws.Range("J15").Value = 8.100363172552354
$ws.Range("L15").Value = 13.10281565819031
$ws.Range("N15").Value = 18.63106951578013
$ws.Range("O15").Value = 29.85592243320005
$ws.Range("B16").Value = 21.79191016566075
$ws.Range("C16").Value = 8.393181476857107
$ws.Range("D16").Value = 13.24593324158358
$ws.Range("E16").Value = 13.19854632267904
$ws.Range("G16").Value = 3.682697996294634
$ws.Range("I16").Value = 27.48699428054194
$ws.Range("J16").Value = 8.104067861107893
$ws.Range("L16").Value = 13.08621278092619
$ws.Range("N16").Value = 18.66880603269407
$ws.Range("O16").Value = 29.87099186072162
$ws.Range("B17").Value = 21.61636090700902
$ws.Range("C17").Value = 8.240942836590797
$ws.Range("D17").Value = 13.24262046665472
$ws.Range("E17").Value = 13.20946510616039
$ws.Range("G17").Value = 3.683683661108416
$ws.Range("I17").Value = 27.52059004861577
$ws.Range("J17").Value = 8.106396737064019
$ws.Range("L17").Value = 13.07652558258356
$ws.Range("N17").Value = 18.69244887681042
$ws.Range("O17").Value = 29.88208792038876
$ws.Range("B18").Value = 21.51503312870346
$ws.Range("C18").Value = 8.151967692329434
$ws.Range("D18").Value = 13.24094753132158
$ws.Range("E18").Value = 13.21590236784785
$ws.Range("G18").Value = 3.6842585216262
$ws.Range("I18").Value = 27.54047154226665
$ws.Range("J18").Value = 8.107756930702017
$ws.Range("L18").Value = 13.07113791455873
$ws.Range("N18").Value = 18.70622903316238
$ws.Range("O18").Value = 29.88914980184867
$ws.Range("B19").Value = 21.4806674737095
$ws.Range("C19").Value = 8.121600100388378
$ws.Range("D19").Value = 13.24042108078879
$ws.Range("E19").Value = 13.21810889896659
$ws.Range("G19").Value = 3.684454523806366
$ws.Range("I19").Value = 27.54729883772887
$ws.Range("J19").Value = 8.108221025855071
$ws.Range("L19").Value = 13.06934547191329
$ws.Range("N19").Value = 18.71092594919817
$ws.Range("O19").Value = 29.89165750590377
$ws.Range("B20").Value = 21.63508612851229
$ws.Range("C20").Value = 8.25729505166262
$ws.Range("D20").Value = 13.24294906722617
$ws.Range("E20").Value = 13.20828652931258
$ws.Range("G20").Value = 3.683577914911623
$ws.Range("I20").Value = 27.51695593938025
$ws.Range("J20").Value = 8.106146684328849
$ws.Range("L20").Value = 13.07753776646293
$ws.Range("N20").Value = 18.68991328749224
$ws.Range("O20").Value = 29.88083636083844
$ws.Range("B21").Value = 22.14940530423604
$ws.Range("C21").Value = 8.696328861397232
$ws.Range("D21").Value = 13.25421673070357
$ws.Range("E21").Value = 13.17716066938573
$ws.Range("G21").Value = 3.680727259485418
$ws.Range("I21").Value = 27.42168263844841
$ws.Range("J21").Value = 8.099424024691045
$ws.Range("L21").Value = 13.1072552827647
$ws.Range("N21").Value = 18.62147873524533
$ws.Range("O21").Value = 29.85260470146284
$ws.Range("B22").Value = 22.48163208002223
$ws.Range("C22").Value = 8.970644423750993
$ws.Range("D22").Value = 13.26363932384534
$ws.Range("E22").Value = 13.15821482613907
$ws.Range("G22").Value = 3.678934397689232
$ws.Range("I22").Value = 27.36440584187422
$ws.Range("J22").Value = 8.095213657552888
$ws.Range("L22").Value = 13.12829350902024
$ws.Range("N22").Value = 18.57836097583129
$ws.Range("O22").Value = 29.84021974048175
$ws.Range("B23").Value = 22.30470864526026
$ws.Range("C23").Value = 8.825393431275236
$ws.Range("D23").Value = 13.25842256828788
$ws.Range("E23").Value = 13.16819887557371
$ws.Range("G23").Value = 3.679884878558108
$ws.Range("I23").Value = 27.39451801689921
$ws.Range("J23").Value = 8.097444084178234
$ws.Range("L23").Value = 13.11691833191021
$ws.Range("N23").Value = 18.60122693777933
$ws.Range("O23").Value = 29.84627424639272
$ws.Range("B24").Value = 21.62662170014739
$ws.Range("C24").Value = 8.249906734602698
$ws.Range("D24").Value = 13.2427997853106
$ws.Range("E24").Value = 13.20881886578822
$ws.Range("G24").Value = 3.683625697254643
$ws.Range("I24").Value = 27.51859715465958
$ws.Range("J24").Value = 8.106259666838286
$ws.Range("L24").Value = 13.07707959210942
$ws.Range("N24").Value = 18.69105904304323
$ws.Range("O24").Value = 29.88140006472625
$ws.Range("B25").Value = 20.88205662611571
$ws.Range("C25").Value = 7.57524823042884
$ws.Range("D25").Value = 13.23484085115008
$ws.Range("E25").Value = 13.25864023875148
$ws.Range("G25").Value = 3.687964307814114
$ws.Range("I25").Value = 27.61865360537207
$ws.Range("J25").Value = 8.116560140105726
$ws.Range("L25").Value = 13.05230738774878
$ws.Range("N25").Value = 18.79489826545732
$ws.Range("O25").Value = 29.94526283954927
